# daily auto push: 2026-01-28 02:31 UTC
#
# A new daily-snapshot row for 2026/01/28 (03:00 slot already present in
# row 710) needs an additional 08:00 slot. The source feed inserts it as
# a brand-new row 711 ("2026/01/28", "水", 8, 33) and every subsequent
# row (old 711..752, the 2026/12/29 .. 2027/01/05 block) shifts down by
# one (new 712..753).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 711 - this shifts rows 711:752 down to 712:753
# and grows the sheet's used range to row 753 automatically.
$ws.Rows("711").Insert()

# Column A in this sheet stores the date as literal text (e.g.
# "2026/01/28"), not a real date serial. Plain `.Value = "2026/01/28"`
# would get smart-parsed into a date by Excel, so instead copy the
# already-text cell immediately above (A710, same calendar date) down
# into the new row - this clones both the text value and its (default)
# formatting without introducing any new style.
$ws.Range("A710").Copy($ws.Range("A711"))

# Remaining three columns for the new row.
$ws.Range("B711").Value = "水"
$ws.Range("C711").Value = 8
$ws.Range("D711").Value = 33
